$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in row 3
$ws.Range("C3").Value = 65
$ws.Range("E3").Value = 377
$ws.Range("F3").Value = 25

# Add formulas for the new row 4, copying the number format used by row 3 (L:P)
$ws.Range("L3:P3").Copy()
$ws.Range("L4:P4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("L4:P4").ClearContents()

$excel.CutCopyMode = $false
